$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column H ("Industries") was set to 1 for rows 36-176; update to 0.
$ws.Range("H36:H176").Value = 0
